$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header cell (E1) to the new header cell (F1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:42:34.312538",
    "2021-10-05 13:42:34.312549",
    "2021-10-05 13:42:34.312553",
    "2021-10-05 13:42:34.312556",
    "2021-10-05 13:42:34.312559",
    "2021-10-05 13:42:34.312562",
    "2021-10-05 13:42:34.312564",
    "2021-10-05 13:42:34.312567",
    "2021-10-05 13:42:34.312570",
    "2021-10-05 13:42:34.312573",
    "2021-10-05 13:42:34.312575",
    "2021-10-05 13:42:34.312578",
    "2021-10-05 13:42:34.312581",
    "2021-10-05 13:42:34.312583",
    "2021-10-05 13:42:34.312586",
    "2021-10-05 13:42:34.312589",
    "2021-10-05 13:42:34.312592",
    "2021-10-05 13:42:34.312594",
    "2021-10-05 13:42:34.312597",
    "2021-10-05 13:42:34.312600",
    "2021-10-05 13:42:34.312602",
    "2021-10-05 13:42:34.312605",
    "2021-10-05 13:42:34.312608",
    "2021-10-05 13:42:34.312610",
    "2021-10-05 13:42:34.312614",
    "2021-10-05 13:42:34.312616",
    "2021-10-05 13:42:34.312619",
    "2021-10-05 13:42:34.312622",
    "2021-10-05 13:42:34.312625",
    "2021-10-05 13:42:34.312627",
    "2021-10-05 13:42:34.312630",
    "2021-10-05 13:42:34.312633",
    "2021-10-05 13:42:34.312636",
    "2021-10-05 13:42:34.312639",
    "2021-10-05 13:42:34.312643",
    "2021-10-05 13:42:34.312647",
    "2021-10-05 13:42:34.312650"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
